## DPLKKPS143-001 / DPLKKPS143-002 - Register Split Balance
## Update the "Lanjutkan ke Verifikasi" register row (sheet DPLKKPS143-002):
##  - SNF_IURAN_PRIBADI (Q2) and SNF_IURAN_PERUSAHAAN (R2) amounts change
##  - the free-text PREPARATION note (F2) is updated to match the new amounts
##  - the NO_REGISTER value (N2) changes
##  - a stray formatted-but-empty cell (Q3) is added below the amounts
## Also flip which tab/cell is active/selected on each sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("DPLKKPS143-001")
$ws2 = $wb.Worksheets.Item("DPLKKPS143-002")

# --- Update the register amounts on DPLKKPS143-002 ---
$ws2.Range("Q2").Value = 2000000
$ws2.Range("R2").Value = 1600000

# --- Update the free-text preparation note to mention the new amounts ---
$prepText = "Username : 33028;`n" + `
    "Password : bni1234;`n" + `
    "Role : 10 - Asisten Settlement;`n" + `
    "Keterangan Perubahan : KEP.TRX.445 melakukan Split Iuran;`n" + `
    "Saldo Nominal Final - Saldo Awal Iuran Pribadi : 2.000.000,00;`n" + `
    "Saldo Nominal Final - Saldo Awal Iuran Perusahaan : 1.600.000,00;`n" + `
    "Saldo Nominal Final - Saldo Awal Iuran Sukarela : 0,00;`n" + `
    "Saldo Nominal Final - Saldo Awal Pengalihan Iuran Karyawan : 0,00;`n" + `
    "Saldo Nominal Final - Saldo Awal Pengalihan Iuran Perusahaan : 0,00;`n" + `
    "Status Register : 1 - Lanjutkan Ke Verifikasi;`n" + `
    "Keterangan Register : KEP.TRX.445 Lanjutkan Verifikasi"
$ws2.Range("F2").Value = $prepText

# --- Update the register number ---
$ws2.Range("N2").Value = "M11220800000011"

# --- Add the new (empty) formatted cell Q3: number format + left/center alignment ---
$ws2.Range("Q3").NumberFormat = "#,##0.00"
$ws2.Range("Q3").HorizontalAlignment = -4131
$ws2.Range("Q3").VerticalAlignment = -4108

# --- Flip the active sheet / selection on each tab ---
# DPLKKPS143-002 becomes inactive, selection moves from R2 to S2
$ws2.Activate()
$ws2.Range("S2").Select()

# DPLKKPS143-001 becomes the active tab, selection moves from G2 to Y2
$ws1.Activate()
$ws1.Range("Y2").Select()
